# Growth reference data and plot function update
# - UK90 curves: reset horizontal scroll (drop topLeftCell="B1"), keep B3 selected
# - Belgium curves: add age 19 & 20 rows for Male and Female growth-reference table
# - Norway curves: add age 19 row for Male and Female growth-reference table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# UK90 curves: scroll back so column A is visible again (removes topLeftCell)
# ---------------------------------------------------------------------------
$wsUK = $wb.Worksheets.Item("UK90 curves")
$wsUK.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1

# ---------------------------------------------------------------------------
# Belgium curves: extend the Male (rows 2-13) and Female (rows 14-25) LMS
# tables with age 19 & 20 entries
# ---------------------------------------------------------------------------
$wsBel = $wb.Worksheets.Item("Belgium curves")
$wsBel.Activate()

# Make room for the two new Male rows (19 & 20) right after the existing
# Male block (ages 7-18) -- this pushes the whole Female block down by 2.
$wsBel.Rows("14:15").Insert()

function Set-Row14Cols {
    param($ws, $row, $gender, $age, $values)
    $ws.Range("A$row").Value2 = $gender
    $ws.Range("B$row").Value2 = $age
    $ws.Range("C$row").Value2 = 1
    $cols = @("D","E","F","G","H","I","J","K","L","M","N")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value2 = $values[$i]
    }
}

# New Male rows (age 19, age 20)
Set-Row14Cols $wsBel 14 "Male" 19 @(180.2, 0.038, 167.8, 169.3, 171.6, 175.6, 180.2, 184.9, 189.2, 191.8, 193.6)
Set-Row14Cols $wsBel 15 "Male" 20 @(180.8, 0.037, 168.6, 170.1, 172.4, 176.3, 180.8, 185.4, 189.6, 192.1, 193.8)

# New Female rows (age 19, age 20) appended after the (now shifted) Female
# block, which runs through row 27.
Set-Row14Cols $wsBel 28 "Female" 19 @(166.4, 0.036, 155.5, 156.8, 158.9, 162.4, 166.4, 170.5, 174.3, 176.6, 178.1)
Set-Row14Cols $wsBel 29 "Female" 20 @(166.5, 0.036, 155.6, 156.9, 159,   162.5, 166.5, 170.6, 174.4, 176.7, 178.2)

$wsBel.Range("G25").Select()

# ---------------------------------------------------------------------------
# Norway curves: extend the Male (rows 2-13) and Female (rows 14-25) LMS
# tables with an age 19 entry
# ---------------------------------------------------------------------------
$wsNor = $wb.Worksheets.Item("Norway curves")
$wsNor.Activate()

# Make room for the one new Male row (19) right after the existing Male
# block (ages 7-18) -- this pushes the whole Female block down by 1.
$wsNor.Rows("14:14").Insert()

Set-Row14Cols $wsNor 14 "Male" 19 @(181, 0.0363, 169.1, 170.5, 172.8, 176.6, 181, 185.5, 189.6, 192.1, 193.8)

# New Female row (age 19) appended after the (now shifted) Female block,
# which runs through row 26.
Set-Row14Cols $wsNor 27 "Female" 19 @(167.2, 0.0366, 156.1, 157.4, 159.5, 163.1, 167.2, 171.4, 175.2, 177.6, 179.1)

$wsNor.Range("H22").Select()

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Turkey curves").Activate()
